$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:43:45"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = "02:43:45"
$ws1.Range("B6").Value = "02:43"
$ws1.Range("D6").Value = 0

$ws1.Range("A7").Value = "02:43:45"
$ws1.Range("B7").Value = "03:54"
$ws1.Range("D7").Value = 71

$ws1.Range("A8").Value = "02:43:45"
$ws1.Range("D8").Value = 78

$ws1.Range("A9").Value = "02:43:45"
$ws1.Range("B9").Value = "04:29"
$ws1.Range("C9").Value = "215_ALUAR"
$ws1.Range("D9").Value = 106
$ws1.Range("E9").Value = "LP1912"

# ----- Sheet 2: LP1912-215 -----
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:43:45"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = "02:43:45"
$ws2.Range("B6").Value = "02:43"
$ws2.Range("D6").Value = 0

$ws2.Range("A7").Value = "02:43:45"
$ws2.Range("B7").Value = "04:29"
$ws2.Range("C7").Value = "215_ALUAR"
$ws2.Range("D7").Value = 106
$ws2.Range("E7").Value = "LP1912"

# ----- Sheet 3: 6203-6173 -----
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:43:45"
